# Replicates the "Made sheet2 the same as reader/sheet2" edit:
# sheet2 gains a new AA column (values 100..129 next to the existing
# D/K columns), its dimension/selection/scroll position move to match,
# and it becomes the active sheet/tab of the workbook (previously it was
# Sheet4 that was the selected tab).

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# --- sheet2: add column AA with values 100..129 on rows 1..30 -------------
$aaValues = @(100,101,102,103,104,105,106,107,108,109,110,111,112,113,114,115,116,117,118,119,120,121,122,123,124,125,126,127,128,129)
for ($i = 0; $i -lt $aaValues.Length; $i++) {
    $ws2.Cells.Item($i + 1, 27).Value = $aaValues[$i]
}

# --- sheet2 becomes the active/selected sheet (was Sheet4 before) ---------
$ws2.Activate()

# --- sheet2's own view: selection is the new AA1:AA30 range ---------------
$ws2.Range("AA1:AA30").Select()

# Scroll the view so column O is the left-most visible column (best effort;
# mirrors the author scrolling right after filling in the new column).
$excel.ActiveWindow.ScrollColumn = 15
$excel.ActiveWindow.ScrollRow = 1
